$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Matt Showman" borrowed CB11 -> now "Johnny Tsunami"
$ws.Range("C2").Value = "Johnny Tsunami"

# Update the active selection to reflect the edited cell
$ws.Range("C3").Select()
